$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project")

# Add the Kanban model link in C5 (B5 already contains "Kanban model")
$ws.Range("C5").Value = "https://www.canva.com/design/DAFdcfvUnxU/WSismOO5mmrhMfZfhGy6sA/edit?utm_content=DAFdcfvUnxU&utm_campaign=designshare&utm_medium=link2&utm_source=sharebutton"

# Move the active selection to C5
$ws.Range("C5").Select()

$wb.Save()
